# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#   - Status (col C) flips from "Ready for handoff" to the handed-back state.
#   - Two new columns get populated for the handed-back rows:
#       F = Latest Target File    (same file that was handed off)
#       G = Latest Handback File  (the localized file that came back)
#     both rendered as live hyperlinks, matching the look of the existing
#     A/D hyperlink cells.
#   - H = Latest Handback DateTime gets a real timestamp instead of the
#     zero-date placeholder.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$localeInfo = @{
    "zh-cn" = @{ HandbackDate = "2016-03-25 10:34:44" };
    "de-de" = @{ HandbackDate = "2016-03-25 10:34:59" };
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $localeInfo[$sheetName]

    # Source/handoff info lives on row 2 - both report rows point at the
    # same just-completed handback, so F/G are identical on rows 2 and 3.
    $sourceFileName = $ws.Range("A2").Value()
    $sourceFileUrl = $ws.Hyperlinks.Item(1).Address

    $handoffFileName = $ws.Range("D2").Value()
    $handoffFileUrl = $ws.Hyperlinks.Item(2).Address
    $handbackFileUrl = $handoffFileUrl.Replace("olhandoff", "olhandback").Replace("ol-handoff", "ol-handback")

    foreach ($row in @(2, 3)) {
        $ws.Range("C$row").Value = $statusText

        $ws.Hyperlinks.Add($ws.Range("F$row"), $sourceFileUrl, "", "", $sourceFileName)
        $ws.Hyperlinks.Add($ws.Range("G$row"), $handbackFileUrl, "", "", $handoffFileName)

        $ws.Range("H$row").Value = $info.HandbackDate
    }
}
